$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177, shifting existing rows 177-224 down to 178-225
$ws.Rows.Item(177).Insert()

# Fill in the new row 177 with the new record's data
$ws.Cells.Item(177, 1).Value = 4
$ws.Cells.Item(177, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(177, 3).Value = "Los Lagos"
$ws.Cells.Item(177, 4).Value = 44722
$ws.Cells.Item(177, 5).Value = 10
$ws.Cells.Item(177, 6).Value = 100112039
$ws.Cells.Item(177, 7).Value = "Ciboulette"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 240
$ws.Cells.Item(177, 11).Value = 2500
$ws.Cells.Item(177, 12).Value = 2500
$ws.Cells.Item(177, 13).Value = 2500
$ws.Cells.Item(177, 14).Value = "$/docena de atados"
$ws.Cells.Item(177, 15).Value = "Región Metropolitana"
$ws.Cells.Item(177, 16).Value = 833
$ws.Cells.Item(177, 17).Value = 3
$ws.Cells.Item(177, 18).Value = "Hortaliza"
